# "changing to online spreadsheet"
# Add two new worker records at the bottom of the worker data table
# (rows 53 and 54), un-hide the previously-hidden "Worker Number"
# helper column, and leave the cursor/selection on the last entry
# typed (H54), matching the author's final editing position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The helper column C was hidden before; the edited workbook shows it.
$ws.Columns.Item(3).Hidden = $false

# New row 53: CESAR VILLARROEL / CONSTANZA ROCHA / 51 / GASFITERÍA / F
$ws.Cells.Item(53, 1).Value = "CESAR VILLARROEL"
$ws.Cells.Item(53, 2).Value = "CONSTANZA ROCHA"
$ws.Cells.Item(53, 3).Value = 51
$ws.Cells.Item(53, 4).Value = "GASFITERÍA"
$ws.Cells.Item(53, 5).Value = "F"

# New row 54: CESAR VILLARROEL / PRUEBA / 52 / GASFITERÍA / M
$ws.Cells.Item(54, 1).Value = "CESAR VILLARROEL"
$ws.Cells.Item(54, 2).Value = "PRUEBA"
$ws.Cells.Item(54, 3).Value = 52
$ws.Cells.Item(54, 4).Value = "GASFITERÍA"
$ws.Cells.Item(54, 5).Value = "M"

# Leave the selection where the author ended up after typing the new data.
$ws.Range("H54").Select()
